$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update URL value (row 2, column B): pythia -> cicada
$ws.Range("B2").Value = "http://fhirfli.dev/fhir/ig/cicada/CodeSystem/VaccineGender"

# Update Date value (row 8, column B)
$ws.Range("B8").Value = "2026-02-11T14:37:07-05:00"

# Insert a new row after row 10 (Contact row) for Jurisdiction
$ws.Rows.Item(11).Insert()

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# Copy the style of the row above (Contact row, style index 2 body style) to the new row
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Re-set the values since paste special (formats) should not affect values, but ensure correctness
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
